# Read.py [Connected, reader active]
#
# A student card was scanned at parking spot #1: mark it Occupied + Registered,
# stamp the "last read" time in E1, and clear out the leftover placeholder
# Student IDs that had been sitting in rows that were never actually
# registered (the reader/ingest script re-writes the whole sheet on each
# scan, so those stale test values get swept back to 0 along the way).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Spot 1 (row 2): occupied + registered, keep its existing student id.
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = "y"

# Stamp the read time (written as a literal timestamp, same as the ingest
# script's original `datetime.now()` write - not a live formula). Go through
# =NOW() and then freeze the computed value so the stored cell is a plain
# number, matching the source file.
$ws.Range("E1").Formula = "=NOW()"
$ws.Range("E1").Value = $ws.Range("E1").Value2

# Row 6's "Occupied" cell gets normalized from a hard literal FALSE to the
# same =FALSE() formula every other still-empty spot uses.
$ws.Range("B6").Formula = "=FALSE()"

# Sweep stale placeholder Student IDs back to 0 for spots that aren't
# actually registered.
$ws.Range("D3:D5").Value = 0
$ws.Range("D7:D12").Value = 0

# Move the active cell to where the reader script left it.
[void]$ws.Range("E4").Select()
